# Weekly update: a new "Fruta / hortaliza" observation is inserted at the
# top of the data block (row 7), pushing the existing rows 7-19 down to
# rows 8-20 and growing the sheet's used range from A1:T19 to A1:T20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; this shifts the previous rows 7..19
# down to 8..20 (carrying their values/formatting with them) and extends
# the sheet dimension automatically.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with this week's record.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44881
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = "Otros"
$ws.Range("I7").Value = 100107002
$ws.Range("J7").Value = "Chirimoya"
$ws.Range("K7").Value = "Cultivar IV Región"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 22000
$ws.Range("O7").Value = 23000
$ws.Range("P7").Value = 22500
$ws.Range("Q7").Value = "`$/caja 12 kilos"
$ws.Range("R7").Value = "Región de Coquimbo"
$ws.Range("S7").Value = 1875
$ws.Range("T7").Value = 12
